# Auto-generated edit script: updates leve profit calculation sheets
# with refreshed market-board price data (per commit: chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 80.666664
$ws.Range("I5").Value = 80.666664
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 80.666664
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 34.333336
$ws.Range("N5").ClearContents()
$ws.Range("H33").Value = 457.31818
$ws.Range("I33").Value = 324.1579
$ws.Range("J33").Value = 1300.6666
$ws.Range("K33").Value = 324.1579
$ws.Range("L33").Value = 1300.6666
$ws.Range("M33").Value = -95.15789999999998
$ws.Range("N33").Value = -1758.6666
$ws.Range("H34").Value = 1769.6
$ws.Range("I34").Value = 1769.6
$ws.Range("K34").Value = 1769.6
$ws.Range("M34").Value = -1566.6
$ws.Range("H36").Value = 1769.6
$ws.Range("I36").Value = 1769.6
$ws.Range("K36").Value = 1769.6
$ws.Range("M36").Value = -1054.6
$ws.Range("H47").Value = 15466.75
$ws.Range("I47").Value = 15466.75
$ws.Range("K47").Value = 15466.75
$ws.Range("M47").Value = -14494.75
$ws.Range("H55").Value = 160.84616
$ws.Range("I55").Value = 60.75
$ws.Range("J55").Value = 321
$ws.Range("K55").Value = 60.75
$ws.Range("L55").Value = 321
$ws.Range("M55").Value = 153.25
$ws.Range("N55").Value = -749
$ws.Range("H58").Value = 804
$ws.Range("I58").Value = 671.1111
$ws.Range("K58").Value = 2013.3333
$ws.Range("M58").Value = -1863.3333
$ws.Range("H64").Value = 7050.3335
$ws.Range("I64").Value = 7140.6
$ws.Range("K64").Value = 7140.6
$ws.Range("M64").Value = -6892.6
$ws.Range("H67").Value = 7050.3335
$ws.Range("I67").Value = 7140.6
$ws.Range("K67").Value = 7140.6
$ws.Range("M67").Value = -6282.6
$ws.Range("H76").Value = 5777.1665
$ws.Range("I76").Value = 4453.5454
$ws.Range("K76").Value = 4453.5454
$ws.Range("M76").Value = -4138.5454
$ws.Range("H79").Value = 5777.1665
$ws.Range("I79").Value = 4453.5454
$ws.Range("K79").Value = 4453.5454
$ws.Range("M79").Value = -3361.5454
$ws.Range("H88").Value = 1000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -1812
$ws.Range("H91").Value = 1000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -3808
$ws.Range("H92").Value = 1350.3636
$ws.Range("J92").Value = 3000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496
$ws.Range("H129").Value = 1799.5
$ws.Range("I129").Value = 1677.1818
$ws.Range("J129").Value = 2068.6
$ws.Range("K129").Value = 5031.5454
$ws.Range("L129").Value = 6205.799999999999
$ws.Range("M129").Value = -31.54539999999997
$ws.Range("N129").Value = -16205.8
$ws.Range("H132").Value = 4329.074
$ws.Range("I132").Value = 4384.0386
$ws.Range("K132").Value = 13152.1158
$ws.Range("M132").Value = -10622.1158
$ws.Range("H135").Value = 1034.2
$ws.Range("I135").Value = 1078.375
$ws.Range("K135").Value = 9705.375
$ws.Range("M135").Value = -7170.375
$ws.Range("H141").Value = 14216.272
$ws.Range("I141").Value = 15327.9
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 45983.7
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = -40803.7
$ws.Range("N141").Value = -19660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1025.5333
$ws.Range("I2").Value = 1025.5333
$ws.Range("K2").Value = 1025.5333
$ws.Range("M2").Value = -912.5333000000001
$ws.Range("H4").Value = 1032.3334
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2232
$ws.Range("H5").Value = 370.27274
$ws.Range("I5").Value = 370.27274
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 370.27274
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -258.27274
$ws.Range("N5").ClearContents()
$ws.Range("H6").Value = 8048.625
$ws.Range("I6").Value = 7769.857
$ws.Range("K6").Value = 7769.857
$ws.Range("M6").Value = -7596.857
$ws.Range("H45").Value = 2045.3334
$ws.Range("I45").Value = 1470.7333
$ws.Range("K45").Value = 1470.7333
$ws.Range("M45").Value = -1093.7333
$ws.Range("H61").Value = 31252154
$ws.Range("I61").Value = 33335498
$ws.Range("K61").Value = 33335498
$ws.Range("M61").Value = -33335286
$ws.Range("H63").Value = 7723.875
$ws.Range("I63").Value = 2291
$ws.Range("J63").Value = 8500
$ws.Range("K63").Value = 2291
$ws.Range("L63").Value = 8500
$ws.Range("M63").Value = -1605
$ws.Range("N63").Value = -9872
$ws.Range("H66").Value = 7723.875
$ws.Range("I66").Value = 2291
$ws.Range("J66").Value = 8500
$ws.Range("K66").Value = 11455
$ws.Range("L66").Value = 42500
$ws.Range("M66").Value = -8023
$ws.Range("N66").Value = -49364
$ws.Range("H74").Value = 250281570
$ws.Range("J74").Value = 3062.5
$ws.Range("L74").Value = 3062.5
$ws.Range("N74").Value = -4810.5
$ws.Range("H77").Value = 250281570
$ws.Range("J77").Value = 3062.5
$ws.Range("L77").Value = 15312.5
$ws.Range("N77").Value = -24048.5
$ws.Range("H88").Value = 26645.5
$ws.Range("I88").Value = 25455.834
$ws.Range("K88").Value = 25455.834
$ws.Range("M88").Value = -25049.834
$ws.Range("H91").Value = 26645.5
$ws.Range("I91").Value = 25455.834
$ws.Range("K91").Value = 25455.834
$ws.Range("M91").Value = -24051.834
$ws.Range("H97").Value = 1255.862
$ws.Range("I97").Value = 1086.4286
$ws.Range("K97").Value = 1086.4286
$ws.Range("M97").Value = -590.4286
$ws.Range("H110").Value = 18615.791
$ws.Range("I110").Value = 24029.412
$ws.Range("J110").Value = 5468.4287
$ws.Range("K110").Value = 24029.412
$ws.Range("L110").Value = 5468.4287
$ws.Range("M110").Value = -21984.412
$ws.Range("N110").Value = -9558.4287
$ws.Range("H116").Value = 1025.5333
$ws.Range("I116").Value = 1025.5333
$ws.Range("K116").Value = 1025.5333
$ws.Range("M116").Value = 1268.4667
$ws.Range("H122").Value = 7941151.5
$ws.Range("I122").Value = 2590
$ws.Range("J122").Value = 9528864
$ws.Range("K122").Value = 7770
$ws.Range("L122").Value = 28586592
$ws.Range("M122").Value = -5320
$ws.Range("N122").Value = -28591492
$ws.Range("H132").Value = 45468410
$ws.Range("I132").Value = 9424.6
$ws.Range("K132").Value = 28273.8
$ws.Range("M132").Value = -25743.8
$ws.Range("H136").Value = 31252154
$ws.Range("I136").Value = 33335498
$ws.Range("K136").Value = 100006494
$ws.Range("M136").Value = -100003944
$ws.Range("H139").Value = 97715
$ws.Range("J139").Value = 97715
$ws.Range("L139").Value = 97715
$ws.Range("N139").Value = -107995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1025.5333
$ws.Range("I3").Value = 1025.5333
$ws.Range("K3").Value = 1025.5333
$ws.Range("M3").Value = -911.5333000000001
$ws.Range("H4").Value = 370.27274
$ws.Range("I4").Value = 370.27274
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 370.27274
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -255.27274
$ws.Range("N4").ClearContents()
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10336
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H81").Value = 22464.715
$ws.Range("J81").Value = 22464.715
$ws.Range("L81").Value = 22464.715
$ws.Range("N81").Value = -24586.715
$ws.Range("H84").Value = 22464.715
$ws.Range("J84").Value = 22464.715
$ws.Range("L84").Value = 67394.145
$ws.Range("N84").Value = -78002.145
$ws.Range("H86").Value = 23033.2
$ws.Range("I86").Value = 51003
$ws.Range("J86").Value = 4386.6665
$ws.Range("K86").Value = 51003
$ws.Range("L86").Value = 4386.6665
$ws.Range("M86").Value = -49880
$ws.Range("N86").Value = -6632.6665
$ws.Range("H89").Value = 23033.2
$ws.Range("I89").Value = 51003
$ws.Range("J89").Value = 4386.6665
$ws.Range("K89").Value = 255015
$ws.Range("L89").Value = 21933.3325
$ws.Range("M89").Value = -249399
$ws.Range("N89").Value = -33165.3325
$ws.Range("H96").Value = 54733
$ws.Range("I96").Value = 9499.5
$ws.Range("J96").Value = 145200
$ws.Range("K96").Value = 9499.5
$ws.Range("L96").Value = 145200
$ws.Range("M96").Value = -6753.5
$ws.Range("N96").Value = -150692
$ws.Range("H107").Value = 2413.1428
$ws.Range("I107").Value = 1973.25
$ws.Range("J107").Value = 2999.6667
$ws.Range("K107").Value = 1973.25
$ws.Range("L107").Value = 2999.6667
$ws.Range("M107").Value = -53.25
$ws.Range("N107").Value = -6839.6667
$ws.Range("H134").Value = 2820.923
$ws.Range("I134").Value = 2820.923
$ws.Range("K134").Value = 8462.769
$ws.Range("M134").Value = -5927.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 286.8125
$ws.Range("J7").Value = 499
$ws.Range("L7").Value = 499
$ws.Range("N7").Value = -725
$ws.Range("H16").Value = 1423.2106
$ws.Range("I16").Value = 1272.8462
$ws.Range("K16").Value = 1272.8462
$ws.Range("M16").Value = -985.8462
$ws.Range("H22").Value = 10407.8
$ws.Range("I22").Value = 20350
$ws.Range("J22").Value = 465.6
$ws.Range("K22").Value = 20350
$ws.Range("L22").Value = 465.6
$ws.Range("M22").Value = -20000
$ws.Range("N22").Value = -1165.6
$ws.Range("H31").Value = 27176208
$ws.Range("I31").Value = 2441.7693
$ws.Range("J31").Value = 37881024
$ws.Range("K31").Value = 2441.7693
$ws.Range("L31").Value = 37881024
$ws.Range("M31").Value = -2146.7693
$ws.Range("N31").Value = -37881614
$ws.Range("H34").Value = 27176208
$ws.Range("I34").Value = 2441.7693
$ws.Range("J34").Value = 37881024
$ws.Range("K34").Value = 2441.7693
$ws.Range("L34").Value = 37881024
$ws.Range("M34").Value = -2239.7693
$ws.Range("N34").Value = -37881428
$ws.Range("H41").Value = 1000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H105").Value = 12664.1
$ws.Range("I105").Value = 2328.75
$ws.Range("K105").Value = 2328.75
$ws.Range("M105").Value = -581.75
$ws.Range("H107").Value = 1211.6923
$ws.Range("I107").Value = 1010.5
$ws.Range("J107").Value = 1533.6
$ws.Range("K107").Value = 1010.5
$ws.Range("L107").Value = 1533.6
$ws.Range("M107").Value = 909.5
$ws.Range("N107").Value = -5373.6
$ws.Range("H113").Value = 1423.2106
$ws.Range("I113").Value = 1272.8462
$ws.Range("K113").Value = 1272.8462
$ws.Range("M113").Value = 897.1538
$ws.Range("H118").Value = 60000
$ws.Range("J118").Value = 60000
$ws.Range("L118").Value = 60000
$ws.Range("N118").Value = -63314
$ws.Range("H122").Value = 2026639.6
$ws.Range("I122").Value = 2346.6365
$ws.Range("K122").Value = 7039.9095
$ws.Range("M122").Value = -4589.9095
$ws.Range("H132").Value = 1807.7354
$ws.Range("I132").Value = 1832.2122
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5496.6366
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -2966.6366
$ws.Range("N132").Value = -8060
$ws.Range("H134").Value = 5200.5
$ws.Range("I134").Value = 4336.2
$ws.Range("J134").Value = 7793.4
$ws.Range("K134").Value = 13008.6
$ws.Range("L134").Value = 23380.2
$ws.Range("M134").Value = -10473.6
$ws.Range("N134").Value = -28450.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5791.2
$ws.Range("J88").Value = 6489
$ws.Range("L88").Value = 19467
$ws.Range("N88").Value = -20323
$ws.Range("H91").Value = 5791.2
$ws.Range("J91").Value = 6489
$ws.Range("L91").Value = 19467
$ws.Range("N91").Value = -22431
$ws.Range("H103").Value = 3007.5
$ws.Range("I103").Value = 1678.6
$ws.Range("J103").Value = 3745.7778
$ws.Range("K103").Value = 5035.799999999999
$ws.Range("L103").Value = 11237.3334
$ws.Range("M103").Value = -4156.799999999999
$ws.Range("N103").Value = -12995.3334
$ws.Range("H113").Value = 1464.1818
$ws.Range("I113").Value = 697
$ws.Range("J113").Value = 1902.5714
$ws.Range("K113").Value = 2091
$ws.Range("L113").Value = 5707.7142
$ws.Range("M113").Value = 79
$ws.Range("N113").Value = -10047.7142
$ws.Range("H116").Value = 4779.8335
$ws.Range("I116").Value = 4559.6665
$ws.Range("K116").Value = 13678.9995
$ws.Range("M116").Value = -10236.9995
$ws.Range("H117").Value = 1134.8
$ws.Range("J117").Value = 1804
$ws.Range("L117").Value = 5412
$ws.Range("N117").Value = -12296
$ws.Range("H124").Value = 2096.6667
$ws.Range("I124").Value = 2233.75
$ws.Range("J124").Value = 1000
$ws.Range("K124").Value = 6701.25
$ws.Range("L124").Value = 3000
$ws.Range("M124").Value = -1791.25
$ws.Range("N124").Value = -12820
$ws.Range("H129").Value = 2094.2856
$ws.Range("J129").Value = 2831.4443
$ws.Range("L129").Value = 8494.332900000001
$ws.Range("N129").Value = -18494.3329
$ws.Range("H130").Value = 2492.4
$ws.Range("I130").Value = 1515
$ws.Range("J130").Value = 3144
$ws.Range("K130").Value = 4545
$ws.Range("L130").Value = 9432
$ws.Range("M130").Value = 475
$ws.Range("N130").Value = -19472
$ws.Range("H132").Value = 3353.9167
$ws.Range("I132").Value = 2186.111
$ws.Range("J132").Value = 4054.6
$ws.Range("K132").Value = 19674.999
$ws.Range("L132").Value = 36491.4
$ws.Range("M132").Value = -17144.999
$ws.Range("N132").Value = -41551.4
$ws.Range("H140").Value = 1811.2
$ws.Range("I140").Value = 1811.2
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 5433.6
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -253.6000000000004
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1017007.56
$ws.Range("I11").Value = 1458399.8
$ws.Range("J11").Value = 23875
$ws.Range("K11").Value = 1458399.8
$ws.Range("L11").Value = 23875
$ws.Range("M11").Value = -1458260.8
$ws.Range("N11").Value = -24153
$ws.Range("H12").Value = 50331.25
$ws.Range("I12").Value = 77163.46000000001
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 77163.46000000001
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -77023.46000000001
$ws.Range("N12").Value = -780
$ws.Range("H64").Value = 29069
$ws.Range("J64").Value = 29069
$ws.Range("L64").Value = 29069
$ws.Range("N64").Value = -29565
$ws.Range("H67").Value = 29069
$ws.Range("J67").Value = 29069
$ws.Range("L67").Value = 29069
$ws.Range("N67").Value = -30785
$ws.Range("H97").Value = 1205.44
$ws.Range("I97").Value = 1126.7826
$ws.Range("J97").Value = 2110
$ws.Range("K97").Value = 1126.7826
$ws.Range("L97").Value = 2110
$ws.Range("M97").Value = -630.7826
$ws.Range("N97").Value = -3102
$ws.Range("H102").Value = 2850.2285
$ws.Range("I102").Value = 2074.3
$ws.Range("K102").Value = 2074.3
$ws.Range("M102").Value = -452.3000000000002
$ws.Range("H113").Value = 4633.222
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H122").Value = 33335888
$ws.Range("I122").Value = 2417.3333
$ws.Range("J122").Value = 83336100
$ws.Range("K122").Value = 7251.999899999999
$ws.Range("L122").Value = 250008300
$ws.Range("M122").Value = -4801.999899999999
$ws.Range("N122").Value = -250013200
$ws.Range("H132").Value = 12149.827
$ws.Range("I132").Value = 8037.9473
$ws.Range("K132").Value = 24113.8419
$ws.Range("M132").Value = -21583.8419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 2949.2727
$ws.Range("I22").Value = 2920.4285
$ws.Range("J22").Value = 2999.75
$ws.Range("K22").Value = 2920.4285
$ws.Range("L22").Value = 2999.75
$ws.Range("M22").Value = -2625.4285
$ws.Range("N22").Value = -3589.75
$ws.Range("H27").Value = 2949.2727
$ws.Range("I27").Value = 2920.4285
$ws.Range("J27").Value = 2999.75
$ws.Range("K27").Value = 2920.4285
$ws.Range("L27").Value = 2999.75
$ws.Range("M27").Value = -2813.4285
$ws.Range("N27").Value = -3213.75
$ws.Range("H46").Value = 1650
$ws.Range("H61").Value = 4137.6
$ws.Range("I61").Value = 2922
$ws.Range("J61").Value = 9000
$ws.Range("K61").Value = 2922
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = -2720
$ws.Range("N61").Value = -9404
$ws.Range("H93").Value = 1686833.1
$ws.Range("I93").Value = 2830.5
$ws.Range("J93").Value = 5054838.5
$ws.Range("K93").Value = 2830.5
$ws.Range("L93").Value = 5054838.5
$ws.Range("M93").Value = -1582.5
$ws.Range("N93").Value = -5057334.5
$ws.Range("H100").Value = 2650.9546
$ws.Range("I100").Value = 2316.2
$ws.Range("J100").Value = 5998.5
$ws.Range("K100").Value = 2316.2
$ws.Range("L100").Value = 5998.5
$ws.Range("M100").Value = -1775.2
$ws.Range("N100").Value = -7080.5
$ws.Range("H113").Value = 4137.6
$ws.Range("I113").Value = 2922
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 2922
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -752
$ws.Range("N113").Value = -13340
$ws.Range("H122").Value = 10421388
$ws.Range("I122").Value = 4708.125
$ws.Range("J122").Value = 31254748
$ws.Range("K122").Value = 14124.375
$ws.Range("L122").Value = 93764244
$ws.Range("M122").Value = -11674.375
$ws.Range("N122").Value = -93769144
$ws.Range("H132").Value = 40003616
$ws.Range("I132").Value = 3766.5715
$ws.Range("K132").Value = 11299.7145
$ws.Range("M132").Value = -8769.7145
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 911403.7
$ws.Range("I136").Value = 1335145.6
$ws.Range("K136").Value = 4005436.8
$ws.Range("M136").Value = -4002886.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H46").Value = 52131
$ws.Range("J46").Value = 52131
$ws.Range("L46").Value = 52131
$ws.Range("N46").Value = -52593
$ws.Range("H62").Value = 11138.889
$ws.Range("I62").Value = 2650
$ws.Range("J62").Value = 21750
$ws.Range("K62").Value = 2650
$ws.Range("L62").Value = 21750
$ws.Range("M62").Value = -2026
$ws.Range("N62").Value = -22998
$ws.Range("H65").Value = 11138.889
$ws.Range("I65").Value = 2650
$ws.Range("J65").Value = 21750
$ws.Range("K65").Value = 13250
$ws.Range("L65").Value = 108750
$ws.Range("M65").Value = -10130
$ws.Range("N65").Value = -114990
$ws.Range("H100").Value = 63126260
$ws.Range("I100").Value = 84167650
$ws.Range("J100").Value = 2099.75
$ws.Range("K100").Value = 168335300
$ws.Range("L100").Value = 4199.5
$ws.Range("M100").Value = -168334759
$ws.Range("N100").Value = -5281.5
$ws.Range("H107").Value = 1634.6
$ws.Range("I107").Value = 1421.6923
$ws.Range("K107").Value = 4265.0769
$ws.Range("M107").Value = -2345.0769
$ws.Range("H113").Value = 953.4167
$ws.Range("I113").Value = 906
$ws.Range("J113").Value = 1095.6666
$ws.Range("K113").Value = 2718
$ws.Range("L113").Value = 3286.9998
$ws.Range("M113").Value = -548
$ws.Range("N113").Value = -7626.9998
$ws.Range("H126").Value = 3703.0454
$ws.Range("I126").Value = 2517.4666
$ws.Range("J126").Value = 6243.5713
$ws.Range("K126").Value = 7552.399800000001
$ws.Range("L126").Value = 18730.7139
$ws.Range("M126").Value = -5082.399800000001
$ws.Range("N126").Value = -23670.7139
$ws.Range("H134").Value = 52131
$ws.Range("J134").Value = 52131
$ws.Range("L134").Value = 156393
$ws.Range("N134").Value = -161463
$ws.Range("H136").Value = 1756.1765
$ws.Range("I136").Value = 1588.5
$ws.Range("K136").Value = 4765.5
$ws.Range("M136").Value = -2215.5
